$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the transistor row: now covers T1 and T2, generic MOSFET description
$ws.Range("A2").Value = "T1,2"
$ws.Range("B2").Value = "N MOSFET"

# New Quantity column
$ws.Range("D1").Value = "Quantity"
$ws.Range("D2").Value = 2

# Update selection to match the new authoring focus
[void]$ws.Range("B2").Select()
